$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1803921568627451
$ws.Cells.Item(2, 3).Value = 0.6039215686274509
$ws.Cells.Item(2, 10).Value = 0.02352941176470588
$ws.Cells.Item(2, 16).Value = 0.1254901960784314
$ws.Cells.Item(2, 19).Value = 0.06666666666666667
$ws.Cells.Item(3, 2).Value = 0.006329113924050633
$ws.Cells.Item(3, 3).Value = 0.0189873417721519
$ws.Cells.Item(3, 10).Value = 0.02531645569620253
$ws.Cells.Item(3, 16).Value = 0.7341772151898734
$ws.Cells.Item(3, 19).Value = 0.2151898734177215
$ws.Cells.Item(4, 10).Value = 0.06451612903225806
$ws.Cells.Item(4, 16).Value = 0.6774193548387096
$ws.Cells.Item(4, 19).Value = 0.2580645161290323
$ws.Cells.Item(6, 2).Value = 0.07352941176470588
$ws.Cells.Item(6, 4).Value = 0.01470588235294118
$ws.Cells.Item(6, 6).Value = 0.0392156862745098
$ws.Cells.Item(6, 10).Value = 0.25
$ws.Cells.Item(6, 15).Value = 0.0196078431372549
$ws.Cells.Item(6, 17).Value = 0.1617647058823529
$ws.Cells.Item(6, 18).Value = 0.04901960784313725
$ws.Cells.Item(6, 19).Value = 0.392156862745098
$ws.Cells.Item(7, 2).Value = 0.07653061224489796
$ws.Cells.Item(7, 4).Value = 0.00510204081632653
$ws.Cells.Item(7, 6).Value = 0.06122448979591837
$ws.Cells.Item(7, 10).Value = 0.1530612244897959
$ws.Cells.Item(7, 15).Value = 0.01530612244897959
$ws.Cells.Item(7, 17).Value = 0.1836734693877551
$ws.Cells.Item(7, 18).Value = 0.07142857142857142
$ws.Cells.Item(7, 19).Value = 0.4336734693877551
$ws.Cells.Item(8, 2).Value = 0.07633587786259542
$ws.Cells.Item(8, 4).Value = 0.01272264631043257
$ws.Cells.Item(8, 6).Value = 0.05089058524173028
$ws.Cells.Item(8, 10).Value = 0.1450381679389313
$ws.Cells.Item(8, 15).Value = 0.02544529262086514
$ws.Cells.Item(8, 17).Value = 0.2239185750636132
$ws.Cells.Item(8, 18).Value = 0.07888040712468193
$ws.Cells.Item(8, 19).Value = 0.3867684478371501
$ws.Cells.Item(9, 2).Value = 0.06077348066298342
$ws.Cells.Item(9, 4).Value = 0.01657458563535912
$ws.Cells.Item(9, 6).Value = 0.06077348066298342
$ws.Cells.Item(9, 10).Value = 0.1104972375690608
$ws.Cells.Item(9, 15).Value = 0.02209944751381215
$ws.Cells.Item(9, 17).Value = 0.1602209944751381
$ws.Cells.Item(9, 18).Value = 0.08839779005524862
$ws.Cells.Item(9, 19).Value = 0.4806629834254144
$ws.Cells.Item(10, 2).Value = 0.09736456808199122
$ws.Cells.Item(10, 4).Value = 0.01464128843338214
$ws.Cells.Item(10, 6).Value = 0.06661786237188873
$ws.Cells.Item(10, 10).Value = 0.1259150805270864
$ws.Cells.Item(10, 15).Value = 0.01317715959004392
$ws.Cells.Item(10, 17).Value = 0.2489019033674963
$ws.Cells.Item(10, 18).Value = 0.05563689604685212
$ws.Cells.Item(10, 19).Value = 0.3777452415812592
$ws.Cells.Item(11, 6).Value = 0.003424657534246575
$ws.Cells.Item(11, 7).Value = 0.1267123287671233
$ws.Cells.Item(11, 10).Value = 0.07876712328767123
$ws.Cells.Item(11, 11).Value = 0.1780821917808219
$ws.Cells.Item(11, 12).Value = 0.5958904109589042
$ws.Cells.Item(11, 19).Value = 0.01712328767123288
$ws.Cells.Item(12, 7).Value = 0.7485714285714286
$ws.Cells.Item(12, 10).Value = 0.2114285714285714
$ws.Cells.Item(12, 12).Value = 0.01142857142857143
$ws.Cells.Item(12, 19).Value = 0.02857142857142857
$ws.Cells.Item(13, 7).Value = 0.6666666666666666
$ws.Cells.Item(13, 10).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.6666666666666666
$ws.Cells.Item(14, 10).Value = 0.3333333333333333
$ws.Cells.Item(15, 6).Value = 0.01762114537444934
$ws.Cells.Item(15, 8).Value = 0.1277533039647577
$ws.Cells.Item(15, 9).Value = 0.05726872246696035
$ws.Cells.Item(15, 10).Value = 0.3700440528634361
$ws.Cells.Item(15, 11).Value = 0.07048458149779736
$ws.Cells.Item(15, 15).Value = 0.08370044052863436
$ws.Cells.Item(15, 19).Value = 0.2731277533039648
$ws.Cells.Item(16, 6).Value = 0.01875
$ws.Cells.Item(16, 8).Value = 0.1625
$ws.Cells.Item(16, 9).Value = 0.1
$ws.Cells.Item(16, 10).Value = 0.45
$ws.Cells.Item(16, 11).Value = 0.1
$ws.Cells.Item(16, 13).Value = 0.03125
$ws.Cells.Item(16, 15).Value = 0.0375
$ws.Cells.Item(16, 19).Value = 0.1
$ws.Cells.Item(17, 6).Value = 0.02476190476190476
$ws.Cells.Item(17, 8).Value = 0.1714285714285714
$ws.Cells.Item(17, 9).Value = 0.08380952380952381
$ws.Cells.Item(17, 10).Value = 0.4419047619047619
$ws.Cells.Item(17, 11).Value = 0.1180952380952381
$ws.Cells.Item(17, 13).Value = 0.01523809523809524
$ws.Cells.Item(17, 15).Value = 0.06666666666666667
$ws.Cells.Item(17, 19).Value = 0.07809523809523809
$ws.Cells.Item(18, 6).Value = 0.01360544217687075
$ws.Cells.Item(18, 8).Value = 0.1768707482993197
$ws.Cells.Item(18, 9).Value = 0.1156462585034014
$ws.Cells.Item(18, 10).Value = 0.4829931972789115
$ws.Cells.Item(18, 11).Value = 0.1020408163265306
$ws.Cells.Item(18, 13).Value = 0.01360544217687075
$ws.Cells.Item(18, 15).Value = 0.04761904761904762
$ws.Cells.Item(18, 19).Value = 0.04761904761904762
$ws.Cells.Item(19, 6).Value = 0.01302931596091205
$ws.Cells.Item(19, 8).Value = 0.1864820846905537
$ws.Cells.Item(19, 9).Value = 0.07573289902280131
$ws.Cells.Item(19, 10).Value = 0.4112377850162867
$ws.Cells.Item(19, 11).Value = 0.1058631921824104
$ws.Cells.Item(19, 13).Value = 0.0252442996742671
$ws.Cells.Item(19, 14).Value = 0.00244299674267101
$ws.Cells.Item(19, 15).Value = 0.07573289902280131
$ws.Cells.Item(19, 19).Value = 0.1042345276872964
